$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1913.7142
$ws.Range("I19").Value = 1819.2
$ws.Range("J19").Value = 2150
$ws.Range("K19").Value = 1819.2
$ws.Range("L19").Value = 2150
$ws.Range("M19").Value = -1644.2
$ws.Range("N19").Value = -2500

$ws.Range("H47").Value = 27000
$ws.Range("J47").Value = 27000
$ws.Range("L47").Value = 27000
$ws.Range("N47").Value = -28944

$ws.Range("H48").Value = 7678.1665
$ws.Range("J48").Value = 7678.1665
$ws.Range("L48").Value = 23034.4995
$ws.Range("N48").Value = -23618.4995

$ws.Range("H56").Value = 7678.1665
$ws.Range("J56").Value = 7678.1665
$ws.Range("L56").Value = 23034.4995
$ws.Range("N56").Value = -24102.4995

$ws.Range("H112").Value = 1717.1428
$ws.Range("J112").Value = 1790.8334
$ws.Range("L112").Value = 5372.5002
$ws.Range("N112").Value = -7588.5002

$ws.Range("H116").Value = 4699.8
$ws.Range("I116").Value = 4642.5713
$ws.Range("K116").Value = 4642.5713
$ws.Range("M116").Value = -1200.5713

$ws.Range("H129").Value = 1123516.9
$ws.Range("J129").Value = 1544736.9
$ws.Range("L129").Value = 4634210.699999999
$ws.Range("N129").Value = -4644210.699999999

$ws.Range("H132").Value = 1401.8088
$ws.Range("I132").Value = 1442.8596
$ws.Range("J132").Value = 1189.091
$ws.Range("K132").Value = 4328.5788
$ws.Range("L132").Value = 3567.273
$ws.Range("M132").Value = -1798.5788
$ws.Range("N132").Value = -8627.272999999999

$ws.Range("H137").Value = 802.91895
$ws.Range("I137").Value = 731.93335
$ws.Range("J137").Value = 851.3182
$ws.Range("K137").Value = 2195.80005
$ws.Range("L137").Value = 2553.9546
$ws.Range("M137").Value = 354.1999500000002
$ws.Range("N137").Value = -7653.9546

$ws.Range("H138").Value = 1602.23
$ws.Range("I138").Value = 752.931
$ws.Range("J138").Value = 2775.0715
$ws.Range("K138").Value = 2258.793
$ws.Range("L138").Value = 8325.2145
$ws.Range("M138").Value = 2881.207
$ws.Range("N138").Value = -18605.2145

$ws.Range("H141").Value = 2424.675
$ws.Range("I141").Value = 815.3871
$ws.Range("K141").Value = 2446.1613
$ws.Range("M141").Value = 2733.8387

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3085.9
$ws.Range("I32").Value = 2619.0854
$ws.Range("J32").Value = 5212.5
$ws.Range("K32").Value = 2619.0854
$ws.Range("L32").Value = 5212.5
$ws.Range("M32").Value = -2332.0854
$ws.Range("N32").Value = -5786.5

$ws.Range("H61").Value = 991.24243
$ws.Range("I61").Value = 787.7917
$ws.Range("J61").Value = 1533.7778
$ws.Range("K61").Value = 787.7917
$ws.Range("L61").Value = 1533.7778
$ws.Range("M61").Value = -575.7917
$ws.Range("N61").Value = -1957.7778

$ws.Range("H74").Value = 921.6
$ws.Range("I74").Value = 906.63416
$ws.Range("K74").Value = 906.63416
$ws.Range("M74").Value = -32.63415999999995

$ws.Range("H77").Value = 921.6
$ws.Range("I77").Value = 906.63416
$ws.Range("K77").Value = 4533.1708
$ws.Range("M77").Value = -165.1707999999999

$ws.Range("H132").Value = 1373.8966
$ws.Range("I132").Value = 1244.1305
$ws.Range("J132").Value = 1871.3334
$ws.Range("K132").Value = 3732.3915
$ws.Range("L132").Value = 5614.0002
$ws.Range("M132").Value = -1202.3915
$ws.Range("N132").Value = -10674.0002

$ws.Range("H136").Value = 991.24243
$ws.Range("I136").Value = 787.7917
$ws.Range("J136").Value = 1533.7778
$ws.Range("K136").Value = 2363.3751
$ws.Range("L136").Value = 4601.3334
$ws.Range("M136").Value = 186.6248999999998
$ws.Range("N136").Value = -9701.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3744.7576
$ws.Range("I105").Value = 3134.4
$ws.Range("J105").Value = 5652.125
$ws.Range("K105").Value = 3134.4
$ws.Range("L105").Value = 5652.125
$ws.Range("M105").Value = -1387.4
$ws.Range("N105").Value = -9146.125

$ws.Range("H107").Value = 8172.95
$ws.Range("I107").Value = 1299.909
$ws.Range("K107").Value = 1299.909
$ws.Range("M107").Value = 620.0909999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2459.7192
$ws.Range("I31").Value = 2315.4324
$ws.Range("J31").Value = 2726.65
$ws.Range("K31").Value = 2315.4324
$ws.Range("L31").Value = 2726.65
$ws.Range("M31").Value = -2020.4324
$ws.Range("N31").Value = -3316.65

$ws.Range("H34").Value = 2459.7192
$ws.Range("I34").Value = 2315.4324
$ws.Range("J34").Value = 2726.65
$ws.Range("K34").Value = 2315.4324
$ws.Range("L34").Value = 2726.65
$ws.Range("M34").Value = -2113.4324
$ws.Range("N34").Value = -3130.65

$ws.Range("H58").Value = 3653.7437
$ws.Range("I58").Value = 1213.5834
$ws.Range("J58").Value = 7558
$ws.Range("K58").Value = 1213.5834
$ws.Range("L58").Value = 7558
$ws.Range("M58").Value = -1010.5834
$ws.Range("N58").Value = -7964

$ws.Range("H132").Value = 1519.1765
$ws.Range("I132").Value = 926.29266
$ws.Range("J132").Value = 2419.4814
$ws.Range("K132").Value = 2778.87798
$ws.Range("L132").Value = 7258.4442
$ws.Range("M132").Value = -248.8779799999998
$ws.Range("N132").Value = -12318.4442

$ws.Range("H134").Value = 1132.7
$ws.Range("I134").Value = 1039.5834
$ws.Range("J134").Value = 1691.4
$ws.Range("K134").Value = 3118.7502
$ws.Range("L134").Value = 5074.200000000001
$ws.Range("M134").Value = -583.7501999999999
$ws.Range("N134").Value = -10144.2

$ws.Range("H136").Value = 3653.7437
$ws.Range("I136").Value = 1213.5834
$ws.Range("J136").Value = 7558
$ws.Range("K136").Value = 3640.7502
$ws.Range("L136").Value = 22674
$ws.Range("M136").Value = -1090.7502
$ws.Range("N136").Value = -27774

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 643.63416
$ws.Range("I5").Value = 618.6061
$ws.Range("J5").Value = 746.875
$ws.Range("K5").Value = 1855.8183
$ws.Range("L5").Value = 2240.625
$ws.Range("M5").Value = -1743.8183
$ws.Range("N5").Value = -2464.625

$ws.Range("H39").Value = 5868.5
$ws.Range("J39").Value = 5868.5
$ws.Range("L39").Value = 17605.5
$ws.Range("N39").Value = -18193.5

$ws.Range("H113").Value = 570.7619
$ws.Range("I113").Value = 604.0476
$ws.Range("J113").Value = 537.4761999999999
$ws.Range("K113").Value = 1812.1428
$ws.Range("L113").Value = 1612.4286
$ws.Range("M113").Value = 357.8571999999999
$ws.Range("N113").Value = -5952.428599999999

$ws.Range("H122").Value = 531.619
$ws.Range("I122").Value = 248.46153
$ws.Range("J122").Value = 658.5517
$ws.Range("K122").Value = 2236.15377
$ws.Range("L122").Value = 5926.9653
$ws.Range("M122").Value = 213.8462300000001
$ws.Range("N122").Value = -10826.9653

$ws.Range("H135").Value = 643.63416
$ws.Range("I135").Value = 618.6061
$ws.Range("J135").Value = 746.875
$ws.Range("K135").Value = 5567.4549
$ws.Range("L135").Value = 6721.875
$ws.Range("M135").Value = -3032.4549
$ws.Range("N135").Value = -11791.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H59").Value = 19615.385
$ws.Range("J59").Value = 20000
$ws.Range("L59").Value = 20000
$ws.Range("N59").Value = -21166

$ws.Range("H109").Value = 10142.5
$ws.Range("J109").Value = 10142.5
$ws.Range("L109").Value = 10142.5
$ws.Range("N109").Value = -12222.5

$ws.Range("H132").Value = 2041.5366
$ws.Range("I132").Value = 1870.1936
$ws.Range("K132").Value = 5610.5808
$ws.Range("M132").Value = -3080.5808

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4995.3794
$ws.Range("I61").Value = 6053
$ws.Range("J61").Value = 1671.4286
$ws.Range("K61").Value = 6053
$ws.Range("L61").Value = 1671.4286
$ws.Range("M61").Value = -5851
$ws.Range("N61").Value = -2075.4286

$ws.Range("H113").Value = 4995.3794
$ws.Range("I113").Value = 6053
$ws.Range("J113").Value = 1671.4286
$ws.Range("K113").Value = 6053
$ws.Range("L113").Value = 1671.4286
$ws.Range("M113").Value = -3883
$ws.Range("N113").Value = -6011.4286

$ws.Range("H132").Value = 1688.0588
$ws.Range("I132").Value = 1608.1915
$ws.Range("K132").Value = 4824.5745
$ws.Range("M132").Value = -2294.5745

$ws.Range("H136").Value = 2337.0476
$ws.Range("I136").Value = 1370
$ws.Range("J136").Value = 4271.143
$ws.Range("K136").Value = 4110
$ws.Range("L136").Value = 12813.429
$ws.Range("M136").Value = -1560
$ws.Range("N136").Value = -17913.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 497.05334
$ws.Range("I132").Value = 479.2903
$ws.Range("J132").Value = 581.7692
$ws.Range("K132").Value = 1437.8709
$ws.Range("L132").Value = 1745.3076
$ws.Range("M132").Value = 1092.1291
$ws.Range("N132").Value = -6805.3076

$ws.Range("H136").Value = 896.725
$ws.Range("I136").Value = 1098.375
$ws.Range("K136").Value = 3295.125
$ws.Range("M136").Value = -745.125
